$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pitch input value (K4): 30 -> 50
$ws.Range("K4").Value = 50

# Update formula for B5: now references B7 instead of B3
$ws.Range("B5").Formula = "=2*TAN(M5/2)*B7"

# Update formula for B6: now references (B7+B13) instead of B4
$ws.Range("B6").Formula = "=2*TAN(M5/2)*(B7+B13)"

# Update formula for B10: divide by (B4-B3) now
$ws.Range("B10").Formula = "=ATAN((B6-B5)/2/(B4-B3))"

# Insert new row 13 content: label "g", formula B13, unit "m"
$ws.Range("A13").Value = "g"
$ws.Range("B13").Formula = "=B3*SIN(PI()/2-M4-M6/2)"
$ws.Range("C13").Value = "m"

# Update the view: zoom to 130%, selection K14
$excel.ActiveWindow.Zoom = 130
$null = $ws.Range("K14").Select()
